$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated K column (column G) values, replacing prior Strike# calc
$kValues = @{
    2 = 2
    3 = 4
    4 = 4
    5 = 4
    6 = 2
    7 = 3
    8 = 5
    9 = 6
    10 = 1
    11 = 0
    12 = 2
    13 = 3
    14 = 0
    15 = 2
    16 = 2
    17 = 0
    18 = 1
    19 = 0
    20 = 2
    21 = 0
    22 = 1
    23 = 0
    24 = 2
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 0
    30 = 0
    31 = 1
    32 = 0
    33 = 1
    34 = 3
    35 = 1
    36 = 0
    37 = 0
    38 = 1
    39 = 1
    40 = 0
    41 = 1
    42 = 0
    43 = 2
    44 = 2
    45 = 0
    46 = 5
    47 = 2
    48 = 3
    49 = 2
    50 = 1
    51 = 1
    52 = 2
    53 = 2
    54 = 2
    55 = 2
    56 = 3
    57 = 1
    58 = 1
    59 = 1
    60 = 1
    61 = 1
    62 = 1
    63 = 2
    64 = 2
    65 = 3
    66 = 2
    67 = 2
    68 = 1
    69 = 1
    70 = 4
    71 = 2
    72 = 2
    73 = 2
    74 = 4
    75 = 1
    76 = 0
    77 = 2
    78 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
